# Add a new row (row 76) of price data to the end of the sheet,
# matching the existing pattern used by the preceding rows (e.g. row 75):
#   - Column A: date stored as text, formatted "YYYY-MM-DD 00:00:00"
#   - Columns B-E: plain numeric values

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(76, 1).Value = "2024-10-17 00:00:00"
$ws.Cells.Item(76, 2).Value = 74700
$ws.Cells.Item(76, 3).Value = 10469.66
$ws.Cells.Item(76, 4).Value = 9265.190000000001
$ws.Cells.Item(76, 5).Value = 7.1214
